$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new data row after the (blank/year-header) row 3, pushing the
#    old "Number of disability persons" row (old row 4) down to row 5, and
#    the source row (old row 5) down to row 6.
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Insert()

# ---------------------------------------------------------------------------
# 1b. Seed the row-label text first (A4/A5) so the shared-string table is
#     built in the same order the source workbook uses, then fill in the
#     title text afterwards.
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("A5").Value = "disabilities Persons "

# ---------------------------------------------------------------------------
# 2. Title row (row 1) - new text, merged across A1:I1
# ---------------------------------------------------------------------------
$ws.Range("A1:I1").Merge()
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Gori Municipality"
$ws.Range("A1:I1").Font.Name = "Arial"
$ws.Range("A1:I1").Font.Bold = $true
$ws.Range("A1:I1").Font.Size = 11
$ws.Range("A1:I1").HorizontalAlignment = -4108
$ws.Range("A1:I1").VerticalAlignment = -4108
$ws.Range("A1:I1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 51

# ---------------------------------------------------------------------------
# 3. Sub-title row (row 2) "(End of year, persons)" - unchanged text
#    but style is refreshed (no fill / no border) and row height reset to
#    default (no explicit custom height).
# ---------------------------------------------------------------------------
$ws.Range("A2").Font.Name = "Arial"
$ws.Range("A2").Font.Size = 10
$ws.Range("A2").Font.ThemeColor = 1
$ws.Range("A2").Interior.Pattern = -4142
$ws.Rows.Item(2).RowHeight = 14.5

# ---------------------------------------------------------------------------
# 4. Row 3 (years header) - A3 blank cell gets the "Sylfaen" font + bottom
#    border; B3:I3 keep their year values/styling (unchanged).
# ---------------------------------------------------------------------------
$ws.Range("A3").Font.Name = "Sylfaen"
$ws.Range("A3").Font.Size = 11
$ws.Range("A3").Font.ThemeColor = 1
$ws.Range("A3").Borders.Item(9).LineStyle = 1
$ws.Range("A3").Borders.Item(9).Weight = 2

# ---------------------------------------------------------------------------
# 5. New row 4: "family with disabilities Persons " + data
# ---------------------------------------------------------------------------
$ws.Range("A4").Font.Name = "Arial"
$ws.Range("A4").Font.Bold = $true
$ws.Range("A4").Font.Size = 11
$ws.Range("A4").Interior.Color = 16777215
$ws.Range("A4").Interior.Pattern = 1
$ws.Range("A4").Borders.Item(7).LineStyle = 1
$ws.Range("A4").Borders.Item(7).Weight = 2
$ws.Range("A4").HorizontalAlignment = -4131
$ws.Range("A4").VerticalAlignment = -4108
$ws.Range("A4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 24.75

$rowB4 = @(2566,2433,2306,2385,2377,2393,2347,2341)
$cols = @("B","C","D","E","F","G","H","I")
for ($i = 0; $i -lt 8; $i++) {
  $addr = $cols[$i] + "4"
  $ws.Range($addr).Value = $rowB4[$i]
  $ws.Range($addr).NumberFormat = "#\ ##0"
  $ws.Range($addr).Font.Name = "Arial"
  $ws.Range($addr).Font.Size = 10
  $ws.Range($addr).Font.Color = 0
  $ws.Range($addr).Interior.Color = 16777215
  $ws.Range($addr).Interior.Pattern = 1
}

# ---------------------------------------------------------------------------
# 6. Row 5 (previously row 4): "disabilities Persons " label + updated data
# ---------------------------------------------------------------------------
$ws.Range("A5").Font.Name = "Arial"
$ws.Range("A5").Font.Bold = $true
$ws.Range("A5").Font.Size = 11
$ws.Range("A5").Interior.Color = 16777215
$ws.Range("A5").Interior.Pattern = 1
$ws.Range("A5").Borders.Item(9).LineStyle = 1
$ws.Range("A5").Borders.Item(9).Weight = 2
$ws.Range("A5").HorizontalAlignment = -4131
$ws.Range("A5").VerticalAlignment = -4108
$ws.Range("A5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 21

$rowB5 = @(2857,2714,2565,2642,2630,2643,2575,2568)
for ($i = 0; $i -lt 8; $i++) {
  $addr = $cols[$i] + "5"
  $ws.Range($addr).Value = $rowB5[$i]
  $ws.Range($addr).NumberFormat = "#\ ##0"
  $ws.Range($addr).Font.Name = "Arial"
  $ws.Range($addr).Font.Size = 10
  $ws.Range($addr).Font.Color = 0
  $ws.Range($addr).Interior.Color = 16777215
  $ws.Range($addr).Interior.Pattern = 1
}
$ws.Range("I5").Borders.Item(9).LineStyle = 1
$ws.Range("I5").Borders.Item(9).Weight = 2

# ---------------------------------------------------------------------------
# 7. Row 6 (previously row 5): Source row, merged A6:H6
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).RowHeight = 27.75

# ---------------------------------------------------------------------------
# 8. Column width / sheet formatting
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.81640625

# ---------------------------------------------------------------------------
# 9. Selection matching the saved view state
# ---------------------------------------------------------------------------
$ws.Range("A1:I1").Select()
